$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 45185 (2023-09-16)
# to 45204 (2023-10-05), keeping the existing date formatting/style.
$ws.Range("C2:C7").Value = 45204
